# Apply the changes described by the diff:
# 1. Metadata sheet: update the "Date" value (B8) to the new timestamp.
# 2. Elements sheet: swap the "Mapping: RIM Mapping" (AK) and
#    "Mapping: Spécification métier vers l'extension ROR OrganizationCreationDate" (AL)
#    columns - both their widths and their cell contents (for every used row).

$wb = $excel.ActiveWorkbook

# --- 1. Update the Date value on the Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value2 = "2024-03-22T16:25:12+00:00"

# --- 2. Swap columns AK and AL on the Elements sheet ---
$elements = $wb.Worksheets.Item("Elements")

# Swap the column widths (AK gets the wider "Spécification métier" column
# width, AL gets the narrower "RIM Mapping" width).
$elements.Columns.Item(37).ColumnWidth = 81.16666667
$elements.Columns.Item(38).ColumnWidth = 24.16666667

# Swap the cell contents for every used row in the sheet.
$lastRow = $elements.UsedRange.Rows.Count
for ($r = 1; $r -le $lastRow; $r++) {
    $akCell = $elements.Cells.Item($r, 37)
    $alCell = $elements.Cells.Item($r, 38)
    $akValue = $akCell.Value2
    $alValue = $alCell.Value2
    $akCell.Value2 = $alValue
    $alCell.Value2 = $akValue
}
